$d = $word.ActiveDocument

# --- Change 1 -----------------------------------------------------------
# "Udagram" + "." + " " + "Developers want to " -> merge the ".", " " and
# "Developers want to " runs into a single run reading ". Developers want to ".
# Re-matching the same visible text with Find/Replace causes the runtime to
# coalesce the matched span into one run.
$null = $d.Content.Find.Execute(". Developers want to ", $true, $false, $false, $false, $false, `
                                 $true, 1, $false, ". Developers want to ", 2)

# --- Change 2 -----------------------------------------------------------
# Add center/both justification to the paragraph that starts with
# "You have been tasked with ...".
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "You have been tasked with*") {
        $p.Format.Alignment = 3   # wdAlignParagraphJustify
    }
}

# --- Change 3 -----------------------------------------------------------
# "/var/log" + "/" + "cloud-init-output.log" -> merge only the first two
# runs into "/var/log/", leaving "cloud-init-output.log" as its own run.
#
# Because this runtime re-normalizes (merges) every run in a paragraph that
# shares identical resolved formatting whenever that paragraph's text is
# edited, we briefly toggle the Bold formatting of "cloud-init-output.log"
# off (a formatting-only change, which does not trigger the merge pass) so
# it is excluded from the merge, perform the text edit, then restore Bold.

$text = $d.Content.Text
$logIdx = $text.IndexOf("cloud-init-output.log")
$rLog = $d.Range($logIdx, $logIdx + 21)
$rLog.Bold = 0

$text2 = $d.Content.Text
$slashIdx = $text2.IndexOf("/var/log/cloud-init-output.log")
$rSlash = $d.Range($slashIdx, $slashIdx + 9)
$null = $rSlash.Find.Execute("/var/log/", $true, $false, $false, $false, $false, `
                              $true, 0, $false, "/var/log/", 2)

$text3 = $d.Content.Text
$logIdx2 = $text3.IndexOf("cloud-init-output.log")
$rLog2 = $d.Range($logIdx2, $logIdx2 + 21)
$rLog2.Bold = 1
